{"js": "// Office.js (Word JavaScript API) script.\n// Applies five bullet-text rewrites to the resume body, matching the\n// author's commit: inline-friendly wording tweaks that keep formatting.\n//\n// Each edit targets a unique, exact source sentence (verified unique in\n// the document) and replaces it in place via Range.insertText(..., \"Replace\"),\n// which preserves the run's existing character formatting (font, size,\n// color, etc.) because it rewrites the text of the same run/range.\n\nconst edits = [\n  {\n    find: \"Managed the full lifecycle of AI-powered solutions, ensuring quality and scalability for a digital platform streamlining end-of-life planning and preserving legacies. \",\n    replace: \"Managed the full lifecycle of AI-powered solutions, collaborating with cross-functional teams to streamline end-of-life planning and preserve legacies. \"\n  },\n  {\n    find: \"Automated deployments using Docker and Kubernetes, and integrated EC2, S3 & CloudWatch for seamless monitoring. \",\n    replace: \"Automated AWS deployments using CodePipeline; integrated EC2, S3 & CloudWatch for seamless monitoring. \"\n  },\n  {\n    find: \"Led a cross-functional team of 5 in sprint planning, progress tracking, and risk management, increasing sprint velocity by 20%. \",\n    replace: \"Led a 5-member team in sprint planning and risk management, increasing sprint velocity by 20% and ensuring quality and scalability of AI solutions. \"\n  },\n  {\n    find: \"Built a GenAI-driven travel app that personalizes itineraries and guides users through contextual audio. \",\n    replace: \"Managed a GenAI-driven travel app project, collaborating with cross-functional teams to personalize itineraries and guide users through contextual audio, while communicating updates to stakeholders. \"\n  },\n  {\n    find: \"Developed the backend in FastAPI, integrating retrieval pipelines with Airflow to ensure low-latency responses for real-time user queries. \",\n    replace: \"Developed the backend in FastAPI, integrating retrieval pipelines and ensuring low-latency responses for real-time user queries. \"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of edits) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + find);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies five bullet-text rewrites to the resume body, matching the\n# author's commit: inline-friendly wording tweaks that keep formatting.\n#\n# Each edit uses Find/Replace (wdReplaceOne) against Document.Content so the\n# existing run formatting is preserved (Find.Execute rewrites the text of\n# the matched range in place rather than deleting/inserting a new run).\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n\n$edits = @(\n  @{\n    Find    = \"Managed the full lifecycle of AI-powered solutions, ensuring quality and scalability for a digital platform streamlining end-of-life planning and preserving legacies. \"\n    Replace = \"Managed the full lifecycle of AI-powered solutions, collaborating with cross-functional teams to streamline end-of-life planning and preserve legacies. \"\n  },\n  @{\n    Find    = \"Automated deployments using Docker and Kubernetes, and integrated EC2, S3 & CloudWatch for seamless monitoring. \"\n    Replace = \"Automated AWS deployments using CodePipeline; integrated EC2, S3 & CloudWatch for seamless monitoring. \"\n  },\n  @{\n    Find    = \"Led a cross-functional team of 5 in sprint planning, progress tracking, and risk management, increasing sprint velocity by 20%. \"\n    Replace = \"Led a 5-member team in sprint planning and risk management, increasing sprint velocity by 20% and ensuring quality and scalability of AI solutions. \"\n  },\n  @{\n    Find    = \"Built a GenAI-driven travel app that personalizes itineraries and guides users through contextual audio. \"\n    Replace = \"Managed a GenAI-driven travel app project, collaborating with cross-functional teams to personalize itineraries and guide users through contextual audio, while communicating updates to stakeholders. \"\n  },\n  @{\n    Find    = \"Developed the backend in FastAPI, integrating retrieval pipelines with Airflow to ensure low-latency responses for real-time user queries. \"\n    Replace = \"Developed the backend in FastAPI, integrating retrieval pipelines and ensuring low-latency responses for real-time user queries. \"\n  }\n)\n\nforeach ($edit in $edits) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $edit.Find\n    $find.Replacement.Text = $edit.Replace\n\n    $found = $find.Execute(\n        $edit.Find,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $edit.Replace,\n        $wdReplaceOne\n    )\n\n    if (-not $found) {\n        throw \"Could not find target text: $($edit.Find)\"\n    }\n}\n"}
